$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:F, rows 2-10 (row 11 is unchanged)
$values = @(
    @(-0.08368629757078778, 0.3438442104993772,  0.4242218510161214,  0.6513231540611168, 0.6522263476946283),
    @( 0.205312877006492,   0.4985604499272415,  0.4677225124169028,  0.6839024143961643, 0.6588478714022412),
    @( 0.0422407319769525,  0.6156476625237515,  0.8090695585645494,  0.8994829395628077, 0.9076125213463587),
    @( 0.1726989577612019,  0.4859904821377091,  0.5789570236283166,  0.7608922549404197, 0.7487137941609299),
    @(-0.05426476430528045, 0.5519219891770294,  0.6681402045837443,  0.8173984368615739, 0.8242260833855453),
    @( 0.0383632470686223,  0.443900067936677,   0.490699680860556,   0.7004995937618779, 0.7085918159736068),
    @(-0.0001125680424399982, 0.5144441553180874, 0.5897234012367705, 0.7679345032206656, 0.7782428017825981),
    @( 0.02398440652176793, 0.5322186494472522,  0.7385541288137565,  0.8593917202380743, 0.880271447287956),
    @(-0.2556182063508592,  0.3897054822362743,  0.4087892466617125,  0.6393662852088093, 0.6081674939718213)
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = 2 + $j   # column B = 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
